$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# ---- Row 2 : header / bucket boundaries + stray "test" note ----
$ws.Range("B2").Value = -1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 99
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 501
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 10001
$ws.Range("J2").Value = '"test"'

# ---- Row 3 : bus type per column ----
$ws.Range("C3").Value = "S"
$ws.Range("D3").Value = "S"
$ws.Range("E3").Value = "L"
$ws.Range("F3").Value = "L"
$ws.Range("G3").Value = "M"
$ws.Range("H3").Value = "M"

# ---- Row 4 : D1 values ----
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 99
$ws.Range("E4").Value = 99
$ws.Range("F4").Value = 99
$ws.Range("G4").Value = 99
$ws.Range("H4").Value = 99

# ---- Row 5 : D2 values ----
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 401
$ws.Range("G5").Value = 401
$ws.Range("H5").Value = 401

# ---- Row 6 : D3 values ----
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 9500

# ---- Row 8 : Price outputs (formulas) ----
$ws.Range("C8").Value = 2500
$ws.Range("D8").Formula = "=F18+D4*H13"
$ws.Range("E8").Formula = "=F18*C18+(E4*H13+E5*H14)*C22"
$ws.Range("F8").Formula = "=F18*C18+(F4*H13+F5*H14)*C22"
$ws.Range("G8").Formula = "=F18*C19+(G4*H13+G5*H14+G6*H15)*C23"
$ws.Range("H8").Formula = "=F18*C19+(H4*H13+H5*H14+H6*H15)*C23"

# ---- Row 9 : Error markers ----
$ws.Range("B9").Value = "X"
$ws.Range("I9").Value = "X"
$ws.Range("J9").Value = "X"

# ---- Row 12 : headings for the lookup tables on the right ----
$ws.Range("F12").Value = "Distance (km):"
$ws.Range("H12").Value = "Kilometer fee (kr/km):"

# ---- Row 15 : close the bracket on the distance bucket label ----
$ws.Range("G15").Value = "[501 - 10000]"

# ---- Row 17 : rename "Initial fee:" -> "Initial fee (IF):" ----
$ws.Range("F17").Value = "Initial fee (IF):"

# ---- Row 25 : spaced-out total price formula description ----
$ws.Range("A25").Value = "Total price = IF * InitialFee + KF * (10 * D1 + 8 * D2 + 6 * D3)"

# ---- Column G : widen to fit the new "[501 - 10000]" label ----
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666

# ---- Selection, matching the saved view in the edited workbook ----
$null = $ws.Range("J29").Select()
